# Saldo.xlsx — reorder the balance block that starts at row 15.
#
# Before (rows 15-32, 18 rows):
#   005646524 EVANGELINA 4422.48   <- moves to the END of the block, value becomes 422.48
#   004313254 GUSTAVO    4292
#   005654122 ELIANE     2369.27
#   004477812 DIEGO      1123.47   <- removed entirely
#   004504449 KELMA      1000
#   ... (rest of the block, unchanged order) ...
#   005044389 CLAUDIA    446.83
#
# After (rows 15-31, 17 rows):
#   004313254 GUSTAVO    4292
#   005654122 ELIANE     2369.27
#   004504449 KELMA      1000
#   ... (rest of the block, unchanged order) ...
#   005044389 CLAUDIA    446.83
#   005646524 EVANGELINA 422.48

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new row order/content for what is currently rows 15-32. This both
# drops the "004477812 / DIEGO / 1123.47" row and moves the
# "005646524 / EVANGELINA" row to the end of the block (with its corrected
# balance of 422.48 instead of 4422.48).
$newBlock = @(
    @("004313254", "GUSTAVO", 4292),
    @("005654122", "ELIANE", 2369.27),
    @("004504449", "KELMA", 1000),
    @("004392159", "RODRIGO", 900.21),
    @("004448303", "NASSIM", 871.71),
    @("004211368", "ILTON", 826.8),
    @("004361159", "HFR", 804.61),
    @("004488571", "CARLOS", 796.03),
    @("004693308", "LAURA", 764.86),
    @("004975924", "SERGIO", 607.4),
    @("004237325", "RICARDO", 606.19),
    @("005591536", "GUSTAVO", 502.81),
    @("004386464", "CARLOS", 498.26),
    @("004862672", "RENATO", 486.97),
    @("004384167", "DOUGLAS", 464.1),
    @("005044389", "CLAUDIA", 446.83),
    @("005646524", "EVANGELINA", 422.48)
)

$startRow = 15

# Write the 17 new rows into what used to be the 18-row block (rows 15-31).
$r = $startRow
foreach ($entry in $newBlock) {
    $accountCell = $ws.Cells.Item($r, 1)
    # Force the account number to stay text (it has a significant leading
    # zero) without leaving a lasting number-format/quote-prefix style on
    # the cell.
    $accountCell.Value = "'" + $entry[0]
    $accountCell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $r = $r + 1
}

# The block shrank by one row (the DIEGO row was removed), so the old last
# row of the block (row 32, now a stale duplicate of row 31's data) needs
# to go away and everything below shifts up.
$ws.Rows($startRow + $newBlock.Count).Delete()
